$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.478.32"
$ws.Range("E2").Value = "  +0.58%  "
$ws.Range("D3").Value = "1.571.39"
$ws.Range("E3").Value = "  +0.41%  "
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("E5").Value = "  -0.10%  "
$ws.Range("D6").Value = "291.22"
$ws.Range("E6").Value = "  +0.18%  "
$ws.Range("D7").Value = "0.3707"
$ws.Range("E7").Value = "  -1.41%  "
$ws.Range("D8").Value = "49.97"
$ws.Range("E8").Value = "  +1.65%  "
$ws.Range("D9").Value = "0.3374"
$ws.Range("E9").Value = "  -0.63%  "
$ws.Range("B10").Value = "Polygon"
$ws.Range("C10").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D10").Value = "1.139"
$ws.Range("E10").Value = "  +0.56%  "
$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").Value = "0.07528"
$ws.Range("E11").Value = "  -0.51%  "
$ws.Range("E12").Value = "  -0.14%  "
$ws.Range("D13").Value = "21.10"
$ws.Range("E13").Value = "  +0.63%  "
$ws.Range("D14").Value = "6.011"
$ws.Range("E14").Value = "  +0.99%  "
$ws.Range("D15").Value = "6.940"
$ws.Range("D16").Value = "1.570.97"
$ws.Range("E16").Value = "  +0.33%  "
$ws.Range("E17").Value = "  -0.90%  "
$ws.Range("D18").Value = "90.47"
$ws.Range("E18").Value = "  +0.85%  "
$ws.Range("D19").Value = "0.06762"
$ws.Range("E19").Value = "  +0.28%  "
$ws.Range("E20").Value = "  -0.14%  "
$ws.Range("D21").Value = "6.314"
$ws.Range("E21").Value = "  +2.06%  "
$ws.Range("D22").Value = "16.41"
$ws.Range("E22").Value = "  -0.83%  "
$ws.Range("D23").Value = "12.21"
$ws.Range("E23").Value = "  +2.48%  "
$ws.Range("D24").Value = "22.476.48"
$ws.Range("E24").Value = "  +0.55%  "
$ws.Range("D25").Value = "2.382"
$ws.Range("E25").Value = "  +0.20%  "
$ws.Range("D26").Value = "2.597"
$ws.Range("E26").Value = "  -3.64%  "
$ws.Range("D27").Value = "20.06"
$ws.Range("E27").Value = "  -0.40%  "
$ws.Range("D28").Value = "149.08"
$ws.Range("E28").Value = "  +1.44%  "
$ws.Range("D29").Value = "5.049"
$ws.Range("E29").Value = "  +0.65%  "
$ws.Range("D30").Value = "125.17"
$ws.Range("E30").Value = "  -0.40%  "
$ws.Range("D31").Value = "1.747.43"
$ws.Range("E31").Value = "  +0.41%  "
$ws.Range("E32").Value = "  +8.59%  "
$ws.Range("D33").Value = "6.209"
$ws.Range("E33").Value = "  +3.19%  "
$ws.Range("D34").Value = "2.012"
$ws.Range("E34").Value = "  -0.18%  "
$ws.Range("D35").Value = "9.744"
$ws.Range("E35").Value = "  -2.95%  "
$ws.Range("E36").Value = "  -1.56%  "
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").Value = "1.372"
$ws.Range("E37").Value = "  -3.77%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "0.02476"
$ws.Range("E38").Value = "  -0.26%  "
$ws.Range("D39").Value = "0.2299"
$ws.Range("E39").Value = "  +0.49%  "
$ws.Range("D40").Value = "0.06518"
$ws.Range("E40").Value = "  +1.22%  "
$ws.Range("D41").Value = "5.430"
$ws.Range("E41").Value = "  +0.77%  "
$ws.Range("D42").Value = "11.27"
$ws.Range("E42").Value = "  +0.36%  "
$ws.Range("D43").Value = "0.6203"
$ws.Range("E43").Value = "  -1.27%  "
$ws.Range("B44").Value = "Frax"
$ws.Range("C44").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D44").Value = "1.002"
$ws.Range("E44").Value = "  -0.12%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "14.01"
$ws.Range("E45").Value = "  +1.07%  "
$ws.Range("D46").Value = "3.807"
$ws.Range("E46").Value = "  +0.20%  "
$ws.Range("D47").Value = "0.5842"
$ws.Range("E47").Value = "  -1.06%  "
$ws.Range("D48").Value = "129.08"
$ws.Range("E48").Value = "  +3.94%  "
$ws.Range("D49").Value = "2.066"
$ws.Range("E49").Value = "  +0.10%  "
$ws.Range("E50").Value = "  -2.61%  "
$ws.Range("D51").Value = "0.07320"
$ws.Range("E51").Value = "  -0.02%  "
